$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers for the two new columns
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header formatting (bold, border, centered) from H1 onto I1:J1
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-15
$data = @(
    @(2, 1, 2),
    @(3, 1, 6),
    @(4, 1, 6),
    @(5, 1, 6),
    @(6, 2, 6),
    @(7, 1, 3),
    @(8, 1, 6),
    @(9, 6, 9),
    @(10, 1, 6),
    @(11, 8, 8),
    @(12, 3, 8),
    @(13, 3, 7),
    @(14, 1, 6),
    @(15, 1, 5)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
